$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Slide 5 ("There are general 'data agnostic' repositories") - add speaker
#    notes describing Dryad / Zenodo / Figshare / Dataverse.
# ---------------------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$notes5 = $s5.NotesPage
$notesBody5 = $notes5.Shapes.AddPlaceholder(2)

$notesText5 = "Dryad is an international open-access repository of research data.  It is a nonprofit organization that provides long-term access to its contents at no cost to users. The base DPC per data submission is `$120 USD. Access is free.`n`nZenodo built and operated by CERN and OpenAIRE to ensure that everyone can join in Open Science.`n`nFigshare is an online open access repository where researchers can preserve and share their research outputs, including figures, datasets, images, and videos. It is free to upload content and free to access, in adherence to the principle of open data. Figshare is one of a number of portfolio businesses supported by Digital Science, a subsidiary of Springer Nature.`n`nDataverse is funded by Harvard with additional support from the Alfred P. Sloan Foundation, National Science Foundation, National Institutes of Health, Helmsley Charitable Trust, IQSS's Henry A. Murray Research Archive, and many others. "

$notesBody5.TextFrame.TextRange.Text = $notesText5

# ---------------------------------------------------------------------------
# 2) Slide 13 (Zenodo pros/cons) - move/resize the "TextBox 4" shape and
#    tweak its wording.
# ---------------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$box = $s13.Shapes.Item("TextBox 4")

$box.Left = 66.0
$box.Top = 208.64895
$box.Width = 852.575
$box.Height = 94.5141

$tr = $box.TextFrame.TextRange
$fullText = $tr.Text
$needle = "not good for discovery"
$startIdx = $fullText.IndexOf($needle)
$sub = $tr.Characters($startIdx + 1, $needle.Length)
$sub.Text = "not (always) good for discovery"
